# Refresh the cryptos price list: update Price (D) / Volume(1h) (E) columns
# for rows 2..51 as per the latest GitHub Actions scrape.
#
# All D/E cells in the source sheet are stored as plain text (inline
# strings), e.g. "601.77" or "  +1.28%  ", not numbers. For values that
# Excel's auto-type-detection would otherwise coerce into a real number
# (plain "123.45"-style strings with a single decimal point), we force the
# cell to Text format first so the literal string round-trips exactly
# (avoiding float precision drift and a wrong stored type), then restore
# the cell's original (default/"Normal") style so we don't leave a stray
# text format behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.944.77'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '3.736.22'
$ws.Range('E3').Value = '  -3.37%  '
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.69%  '
$ws.Range('D7').Value = '3.734.50'
$ws.Range('E7').Value = '  -3.77%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +2.00%  '
$ws.Range('E10').Value = '  +2.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.33'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.02%  '
$ws.Range('E12').Value = '  -1.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.15'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('D15').Value = '4.359.95'
$ws.Range('E15').Value = '  -2.68%  '
$ws.Range('D16').Value = '3.735.73'
$ws.Range('E16').Value = '  -2.06%  '
$ws.Range('D17').Value = '68.883.13'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.24'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('E19').Value = '  +0.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '497.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.29'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +11.14%  '
$ws.Range('E23').Value = '  -1.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.31'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000139'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.35'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.73%  '
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('E30').Value = '  -0.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.02'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.26%  '
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.79'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.88%  '
$ws.Range('D34').Value = '3.886.43'
$ws.Range('E34').Value = '  -2.12%  '
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('D36').Value = '3.666.93'
$ws.Range('E36').Value = '  -3.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.45%  '
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.82'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.133'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.53%  '
$ws.Range('E41').Value = '  -1.84%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '436.05'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '49.01'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('E44').Value = '  -1.23%  '
$ws.Range('E45').Value = '  -0.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.40'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.59'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '142.18'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.04%  '
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range('D51').Value = '2.745.01'
$ws.Range('E51').Value = '  -3.82%  '
